$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 1633
$ws.Range("I20").Value = 1633
$ws.Range("K20").Value = 1633
$ws.Range("M20").Value = -1403

# Row 29
$ws.Range("H29").Value = 3820.2
$ws.Range("J29").Value = 3820.2
$ws.Range("L29").Value = 11460.6
$ws.Range("N29").Value = -12022.6

# Row 32
$ws.Range("H32").Value = 5558758.5
$ws.Range("I32").Value = 1183.1666
$ws.Range("J32").Value = 8337546
$ws.Range("K32").Value = 1183.1666
$ws.Range("L32").Value = 8337546
$ws.Range("M32").Value = -857.1666
$ws.Range("N32").Value = -8338198

# Row 33
$ws.Range("H33").Value = 350.75
$ws.Range("I33").Value = 346.27274
$ws.Range("K33").Value = 346.27274
$ws.Range("M33").Value = -117.27274

# Row 34
$ws.Range("H34").Value = 2006.1428
$ws.Range("I34").Value = 2006.1428
$ws.Range("K34").Value = 2006.1428
$ws.Range("M34").Value = -1803.1428

# Row 35
$ws.Range("H35").Value = 1633
$ws.Range("I35").Value = 1633
$ws.Range("K35").Value = 1633
$ws.Range("M35").Value = -1254

# Row 36
$ws.Range("H36").Value = 2006.1428
$ws.Range("I36").Value = 2006.1428
$ws.Range("K36").Value = 2006.1428
$ws.Range("M36").Value = -1291.1428

# Row 38
$ws.Range("H38").Value = 3275.5625
$ws.Range("I38").Value = 3600.6428
$ws.Range("K38").Value = 10801.9284
$ws.Range("M38").Value = -10429.9284

# Row 39
$ws.Range("H39").Value = 3829.818
$ws.Range("I39").Value = 5888.4287
$ws.Range("J39").Value = 227.25
$ws.Range("K39").Value = 17665.2861
$ws.Range("L39").Value = 681.75
$ws.Range("M39").Value = -17369.2861
$ws.Range("N39").Value = -1273.75

# Row 40
$ws.Range("H40").Value = 3800.9333
$ws.Range("I40").Value = 3063
$ws.Range("J40").Value = 4365.2354
$ws.Range("K40").Value = 3063
$ws.Range("L40").Value = 4365.2354
$ws.Range("M40").Value = -2888
$ws.Range("N40").Value = -4715.2354

# Row 41
$ws.Range("H41").Value = 1237.8334
$ws.Range("I41").Value = 1015.55554
$ws.Range("K41").Value = 1015.55554
$ws.Range("M41").Value = -575.55554

# Row 42
$ws.Range("H42").Value = 200647.4
$ws.Range("I42").Value = 286053.56
$ws.Range("J42").Value = 1366.3334
$ws.Range("K42").Value = 858160.6799999999
$ws.Range("L42").Value = 4099.0002
$ws.Range("M42").Value = -857930.6799999999
$ws.Range("N42").Value = -4559.0002

# Row 47
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21944

# Row 86
$ws.Range("H86").Value = 5266946.5
$ws.Range("I86").Value = 4722.5
$ws.Range("K86").Value = 4722.5
$ws.Range("M86").Value = -3599.5

# Row 87
$ws.Range("H87").Value = 21874.875
$ws.Range("J87").Value = 21874.875
$ws.Range("L87").Value = 21874.875
$ws.Range("N87").Value = -24370.875

# Row 89
$ws.Range("H89").Value = 5266946.5
$ws.Range("I89").Value = 4722.5
$ws.Range("K89").Value = 23612.5
$ws.Range("M89").Value = -17996.5

# Row 90
$ws.Range("H90").Value = 21874.875
$ws.Range("J90").Value = 21874.875
$ws.Range("L90").Value = 65624.625
$ws.Range("N90").Value = -78104.625

# Row 137
$ws.Range("H137").Value = 3418.6865
$ws.Range("I137").Value = 3641.6345
$ws.Range("J137").Value = 2645.8
$ws.Range("K137").Value = 10924.9035
$ws.Range("L137").Value = 7937.400000000001
$ws.Range("M137").Value = -8374.9035
$ws.Range("N137").Value = -13037.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12991998
$ws.Range("I32").Value = 14290861
$ws.Range("K32").Value = 14290861
$ws.Range("M32").Value = -14290574

# Row 61
$ws.Range("H61").Value = 2942.6875
$ws.Range("I61").Value = 1590.8334
$ws.Range("K61").Value = 1590.8334
$ws.Range("M61").Value = -1378.8334

# Row 74
$ws.Range("H74").Value = 2948.5
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126

# Row 77
$ws.Range("H77").Value = 2948.5
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632

# Row 122
$ws.Range("H122").Value = 4691.7646
$ws.Range("J122").Value = 4928.4287
$ws.Range("L122").Value = 14785.2861
$ws.Range("N122").Value = -19685.2861

# Row 136
$ws.Range("H136").Value = 2942.6875
$ws.Range("I136").Value = 1590.8334
$ws.Range("K136").Value = 4772.5002
$ws.Range("M136").Value = -2222.5002

$ws = $wb.Worksheets.Item("BSM")
# Row 106
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524

# Row 134
$ws.Range("H134").Value = 3163.17
$ws.Range("I134").Value = 2841.175
$ws.Range("K134").Value = 8523.525000000001
$ws.Range("M134").Value = -5988.525000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5797.154
$ws.Range("I16").Value = 6277.4546
$ws.Range("K16").Value = 6277.4546
$ws.Range("M16").Value = -5990.4546

# Row 31
$ws.Range("H31").Value = 2016.9584
$ws.Range("I31").Value = 1905.65
$ws.Range("K31").Value = 1905.65
$ws.Range("M31").Value = -1610.65

# Row 34
$ws.Range("H34").Value = 2016.9584
$ws.Range("I34").Value = 1905.65
$ws.Range("K34").Value = 1905.65
$ws.Range("M34").Value = -1703.65

# Row 58
$ws.Range("H58").Value = 1936.8077
$ws.Range("I58").Value = 1639.875
$ws.Range("K58").Value = 1639.875
$ws.Range("M58").Value = -1436.875

# Row 113
$ws.Range("H113").Value = 5797.154
$ws.Range("I113").Value = 6277.4546
$ws.Range("K113").Value = 6277.4546
$ws.Range("M113").Value = -4107.4546

# Row 122
$ws.Range("H122").Value = 381431.22
$ws.Range("I122").Value = 638852.9
$ws.Range("J122").Value = 6999.727
$ws.Range("K122").Value = 1916558.7
$ws.Range("L122").Value = 20999.181
$ws.Range("M122").Value = -1914108.7
$ws.Range("N122").Value = -25899.181

# Row 132
$ws.Range("H132").Value = 2823.95
$ws.Range("I132").Value = 2977.5625
$ws.Range("J132").Value = 2209.5
$ws.Range("K132").Value = 8932.6875
$ws.Range("L132").Value = 6628.5
$ws.Range("M132").Value = -6402.6875
$ws.Range("N132").Value = -11688.5

# Row 134
$ws.Range("H134").Value = 4359.7036
$ws.Range("I134").Value = 4230.2354
$ws.Range("J134").Value = 4579.8
$ws.Range("K134").Value = 12690.7062
$ws.Range("L134").Value = 13739.4
$ws.Range("M134").Value = -10155.7062
$ws.Range("N134").Value = -18809.4

# Row 136
$ws.Range("H136").Value = 1936.8077
$ws.Range("I136").Value = 1639.875
$ws.Range("K136").Value = 4919.625
$ws.Range("M136").Value = -2369.625

$ws = $wb.Worksheets.Item("CUL")
# Row 95
$ws.Range("H95").Value = 7500
$ws.Range("J95").Value = 7500
$ws.Range("L95").Value = 22500
$ws.Range("N95").Value = -26618

# Row 107
$ws.Range("H107").Value = 367.25
$ws.Range("I107").Value = 345
$ws.Range("J107").Value = 389.5
$ws.Range("K107").Value = 1035
$ws.Range("L107").Value = 1168.5
$ws.Range("M107").Value = 885
$ws.Range("N107").Value = -5008.5

# Row 111
$ws.Range("H111").Value = 2586.25
$ws.Range("I111").Value = 1175
$ws.Range("K111").Value = 3525
$ws.Range("M111").Value = -458

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 700
$ws.Range("I31").Value = 371.42856
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 371.42856
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -79.42856
$ws.Range("N31").Value = -3584

# Row 37
$ws.Range("H37").Value = 700
$ws.Range("I37").Value = 371.42856
$ws.Range("J37").Value = 3000
$ws.Range("K37").Value = 371.42856
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = -94.42856
$ws.Range("N37").Value = -3554

# Row 122
$ws.Range("H122").Value = 4124.1177
$ws.Range("J122").Value = 5363.1875
$ws.Range("L122").Value = 16089.5625
$ws.Range("N122").Value = -20989.5625

# Row 132
$ws.Range("H132").Value = 2776.257
$ws.Range("I132").Value = 2748.9666
$ws.Range("K132").Value = 8246.899800000001
$ws.Range("M132").Value = -5716.899800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3475.121
$ws.Range("I136").Value = 3475.121
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10425.363
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7875.363000000001
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 81400
$ws.Range("J135").Value = 81400
$ws.Range("L135").Value = 81400
$ws.Range("N135").Value = -91540

# Row 136
$ws.Range("H136").Value = 1771.1052
$ws.Range("I136").Value = 967.38464
$ws.Range("K136").Value = 2902.15392
$ws.Range("M136").Value = -352.1539199999997

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
